$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.501.10"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.628.51"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'211.47"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'22.89"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").Value = "1.860.31"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.632.05"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'0.557"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.502.32"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "'228.24"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  +7.07%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "'148.96"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "'6.85"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "'15.56"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "1.463.13"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").Value = "'0.874"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").Value = "'0.917"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "'1.02"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "'67.99"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "1.769.34"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'1.73"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").Value = "'87.35"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0993"
$ws.Range("E51").Value = "  +0.21%  "
